$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.8.1 -> 1.8.2
$ws.Range("B3").Value = "1.8.2"

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Experimental: true -> (cleared)
$ws.Range("B7").ClearContents()

# Date: 2024-01-18 -> 2025-11-18 (write via formula then flatten to a literal so it
# stays a plain text value instead of being auto-converted to a date serial number)
$ws.Range("B8").Formula = '="2025-11-18"'
$ws.Range("B8").Copy()
$ws.Range("B8").PasteSpecial(-4163)
$excel.CutCopyMode = 0
